$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Add files via upload" commit bundles three unrelated row edits on the
# "Export" sheet (Conta / Nome / Saldo):
#   1. Rows 3-4 (RICARDO 3839.41 / CAIO 1516.5) collapse into a single row
#      for DIOGO 2733.11.
#   2. The row that used to hold DIOGO 733.11 (row 10) is rewritten to hold
#      RENAN 794.51.
#   3. The old RENAN row (005216881 / 53.61, row 118) is removed entirely.
#
# Apply the edits from the bottom of the sheet upward so the row numbers
# referenced below stay valid while we work.

# 3) Remove the old RENAN (005216881 / 53.61) row entirely.
$ws.Rows.Item(118).Delete()

# 2) Turn the old DIOGO (733.11) row into the new RENAN (794.51) row.
# "Conta" values are account numbers with significant leading zeros, so the
# column has to be forced to Text before writing them, otherwise Excel
# would store them as plain numbers and drop the leading zero.
$ws.Cells.Item(10, 1).NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "005216881"
$ws.Cells.Item(10, 2).Value = "RENAN"
$ws.Cells.Item(10, 3).Value = 794.51

# 1) Drop the CAIO row, then rewrite what was the RICARDO row as DIOGO.
$ws.Rows.Item(4).Delete()
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004550415"
$ws.Cells.Item(3, 2).Value = "DIOGO"
$ws.Cells.Item(3, 3).Value = 2733.11
